$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.987.15"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "1.750.94"
$ws.Range("E3").Value = "  -3.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3778"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.83%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3356"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.41"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.112"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07212"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.150"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.125"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.63%  "

$ws.Range("D16").Value = "1.754.31"
$ws.Range("E16").Value = "  -3.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001056"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06586"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.250"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.87%  "

$ws.Range("D23").Value = "28.039.79"
$ws.Range("E23").Value = "  -0.59%  "

$ws.Range("E24").Value = "  -6.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.401"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.37%  "

$ws.Range("E27").Value = "  -7.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.325"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.00%  "

$ws.Range("D29").Value = "1.956.54"
$ws.Range("E29").Value = "  -3.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.249"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -15.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.017"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.784"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08640"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.36%  "

$ws.Range("E35").Value = "  -7.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6685"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02316"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06184"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.150"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2110"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.211"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.449"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.988"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.826"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6046"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.017"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.32%  "

$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.179"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.96%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07151"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.87%  "
